# Scheduled data refresh: update Leve profit-tracking sheets with refreshed
# market-board price data (currentAveragePrice / NQ / HQ) and the derived
# LevePriceNQ/HQ and LeveProfitNQ/HQ columns for the affected leve rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5062.5557
$ws.Range("I19").Value = 2332.875
$ws.Range("K19").Value = 2332.875
$ws.Range("M19").Value = -2157.875
$ws.Range("H134").Value = 144400
$ws.Range("J134").Value = 144400
$ws.Range("L134").Value = 144400
$ws.Range("N134").Value = -154540
$ws.Range("H137").Value = 1415.0555
$ws.Range("I137").Value = 1231.4667
$ws.Range("J137").Value = 2333
$ws.Range("K137").Value = 3694.4001
$ws.Range("L137").Value = 6999
$ws.Range("M137").Value = -1144.4001
$ws.Range("N137").Value = -12099
$ws.Range("H138").Value = 7067376
$ws.Range("I138").Value = 12502018
$ws.Range("K138").Value = 37506054
$ws.Range("M138").Value = -37500914
$ws.Range("H139").Value = 99924
$ws.Range("J139").Value = 99924
$ws.Range("L139").Value = 99924
$ws.Range("N139").Value = -110204
$ws.Range("H140").Value = 72950
$ws.Range("J140").Value = 72950
$ws.Range("L140").Value = 72950
$ws.Range("N140").Value = -83310

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1080.5333
$ws.Range("I2").Value = 1214.3
$ws.Range("J2").Value = 813
$ws.Range("K2").Value = 1214.3
$ws.Range("L2").Value = 813
$ws.Range("M2").Value = -1101.3
$ws.Range("N2").Value = -1039
$ws.Range("H80").Value = 34637.332
$ws.Range("J80").Value = 35588.2
$ws.Range("L80").Value = 35588.2
$ws.Range("N80").Value = -37584.2
$ws.Range("H83").Value = 34637.332
$ws.Range("J83").Value = 35588.2
$ws.Range("L83").Value = 106764.6
$ws.Range("N83").Value = -116748.6
$ws.Range("H102").Value = 1421.125
$ws.Range("I102").Value = 1404.7826
$ws.Range("J102").Value = 1797
$ws.Range("K102").Value = 1404.7826
$ws.Range("L102").Value = 1797
$ws.Range("M102").Value = 217.2174
$ws.Range("N102").Value = -5041
$ws.Range("H110").Value = 5588.619
$ws.Range("I110").Value = 5523.3125
$ws.Range("J110").Value = 5797.6
$ws.Range("K110").Value = 5523.3125
$ws.Range("L110").Value = 5797.6
$ws.Range("M110").Value = -3478.3125
$ws.Range("N110").Value = -9887.6
$ws.Range("H116").Value = 1080.5333
$ws.Range("I116").Value = 1214.3
$ws.Range("J116").Value = 813
$ws.Range("K116").Value = 1214.3
$ws.Range("L116").Value = 813
$ws.Range("M116").Value = 1079.7
$ws.Range("N116").Value = -5401

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1080.5333
$ws.Range("I3").Value = 1214.3
$ws.Range("J3").Value = 813
$ws.Range("K3").Value = 1214.3
$ws.Range("L3").Value = 813
$ws.Range("M3").Value = -1100.3
$ws.Range("N3").Value = -1041
$ws.Range("H82").Value = 23452.666
$ws.Range("I82").Value = 3097.2
$ws.Range("K82").Value = 3097.2
$ws.Range("M82").Value = -2714.2
$ws.Range("H85").Value = 23452.666
$ws.Range("I85").Value = 3097.2
$ws.Range("K85").Value = 3097.2
$ws.Range("M85").Value = -1771.2
$ws.Range("H92").Value = 69333.336
$ws.Range("I92").Value = 30000
$ws.Range("K92").Value = 30000
$ws.Range("M92").Value = -27504
$ws.Range("H94").Value = 7143831.5
$ws.Range("I94").Value = 8621452
$ws.Range("J94").Value = 1998.3334
$ws.Range("K94").Value = 8621452
$ws.Range("L94").Value = 1998.3334
$ws.Range("M94").Value = -8621001
$ws.Range("N94").Value = -2900.3334
$ws.Range("H99").Value = 1603.7894
$ws.Range("I99").Value = 1491.2667
$ws.Range("J99").Value = 2025.75
$ws.Range("K99").Value = 1491.2667
$ws.Range("L99").Value = 2025.75
$ws.Range("M99").Value = 6.733300000000099
$ws.Range("N99").Value = -5021.75
$ws.Range("H105").Value = 31251270
$ws.Range("I105").Value = 35715456
$ws.Range("J105").Value = 1975
$ws.Range("K105").Value = 35715456
$ws.Range("L105").Value = 1975
$ws.Range("M105").Value = -35713709
$ws.Range("N105").Value = -5469

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 43422.54
$ws.Range("J64").Value = 43422.54
$ws.Range("L64").Value = 43422.54
$ws.Range("N64").Value = -43918.54
$ws.Range("H67").Value = 43422.54
$ws.Range("J67").Value = 43422.54
$ws.Range("L67").Value = 43422.54
$ws.Range("N67").Value = -45138.54
$ws.Range("H105").Value = 1778.8334
$ws.Range("I105").Value = 1968.25
$ws.Range("K105").Value = 1968.25
$ws.Range("M105").Value = -221.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 81111
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H131").Value = 1384.5555
$ws.Range("J131").Value = 2000
$ws.Range("L131").Value = 6000
$ws.Range("N131").Value = -16080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 43625.5
$ws.Range("J93").Value = 43625.5
$ws.Range("L93").Value = 43625.5
$ws.Range("N93").Value = -47369.5
$ws.Range("H97").Value = 2267.9333
$ws.Range("I97").Value = 2314.111
$ws.Range("J97").Value = 2198.6667
$ws.Range("K97").Value = 2314.111
$ws.Range("L97").Value = 2198.6667
$ws.Range("M97").Value = -1818.111
$ws.Range("N97").Value = -3190.6667
$ws.Range("H113").Value = 15332.333
$ws.Range("I113").Value = 23601.8
$ws.Range("J113").Value = 4995.5
$ws.Range("K113").Value = 23601.8
$ws.Range("L113").Value = 4995.5
$ws.Range("M113").Value = -21431.8
$ws.Range("N113").Value = -9335.5
$ws.Range("H122").Value = 204656.8
$ws.Range("I122").Value = 226285.33
$ws.Range("K122").Value = 678855.99
$ws.Range("M122").Value = -676405.99
$ws.Range("H132").Value = 1252129.5
$ws.Range("J132").Value = 2672.7
$ws.Range("L132").Value = 8018.099999999999
$ws.Range("N132").Value = -13078.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 15218.5
$ws.Range("I61").Value = 17744.182
$ws.Range("K61").Value = 17744.182
$ws.Range("M61").Value = -17542.182
$ws.Range("H82").Value = 829.4286
$ws.Range("I82").Value = 587.1818
$ws.Range("J82").Value = 1095.9
$ws.Range("K82").Value = 587.1818
$ws.Range("L82").Value = 1095.9
$ws.Range("M82").Value = -226.1818
$ws.Range("N82").Value = -1817.9
$ws.Range("H85").Value = 829.4286
$ws.Range("I85").Value = 587.1818
$ws.Range("J85").Value = 1095.9
$ws.Range("K85").Value = 587.1818
$ws.Range("L85").Value = 1095.9
$ws.Range("M85").Value = 660.8182
$ws.Range("N85").Value = -3591.9
$ws.Range("H113").Value = 15218.5
$ws.Range("I113").Value = 17744.182
$ws.Range("K113").Value = 17744.182
$ws.Range("M113").Value = -15574.182
$ws.Range("H136").Value = 5558465.5
$ws.Range("I136").Value = 2618.56
$ws.Range("J136").Value = 18185390
$ws.Range("K136").Value = 7855.68
$ws.Range("L136").Value = 54556170
$ws.Range("M136").Value = -5305.68
$ws.Range("N136").Value = -54561270

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5965.0415
$ws.Range("I107").Value = 1644.1765
$ws.Range("J107").Value = 16458.572
$ws.Range("K107").Value = 4932.529500000001
$ws.Range("L107").Value = 49375.716
$ws.Range("M107").Value = -3012.529500000001
$ws.Range("N107").Value = -53215.716
$ws.Range("H122").Value = 3581.6428
$ws.Range("I122").Value = 3472.5386
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 10417.6158
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -7967.6158
$ws.Range("N122").Value = -19900
$ws.Range("H126").Value = 3794.1667
$ws.Range("I126").Value = 3794.1667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 11382.5001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -8912.500100000001
$ws.Range("N126").ClearContents()
